# The document is a single-column table of numbers. This edit:
#   1. Updates rows 1-3 to "0M" (memory sizes now reported in MB suffix).
#   2. Inserts 10 new per-benchmark-iteration rows right after row 3
#      (a GC-pause-time breakdown that used to be crammed, tab-separated,
#      into two later rows).
#   3. Collapses the two tab-separated rows near the end down to their
#      first value only ("100" and "0" respectively).
#   4. Fills in the trailing empty row with "4532".

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. First three rows: single scalar -> "0M" ---------------------------
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- 2. Insert 10 new rows directly after row 3 ----------------------------
# Table.Rows.Add(beforeRow) inserts immediately *before* beforeRow, and
# repeated calls against the same beforeRow stack in LIFO order, so we walk
# the desired values in reverse to land them in forward order.
$newValues = @("102", "0.00002", "0.00012", "0.00005", "0.00003", "0.00003", "0.00004", "0.00009", "0.00408", "100.0")
$reversedValues = $newValues[($newValues.Length - 1)..0]

$beforeRow = $t.Rows.Item(4)
foreach ($v in $reversedValues) {
    $newRow = $t.Rows.Add($beforeRow)
    $newRow.Cells.Item(1).Range.Text = $v
}

# --- 3 & 4. Last three rows: collapse tab-runs, fill the empty trailer ----
$n = $t.Rows.Count
$t.Rows.Item($n - 2).Cells.Item(1).Range.Text = "100"
$t.Rows.Item($n - 1).Cells.Item(1).Range.Text = "0"
$t.Rows.Item($n).Cells.Item(1).Range.Text = "4532"
